$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E3").Value = 22
$ws.Range("E5").Value = 3
$ws.Range("E17").Value = 99
$ws.Range("E19").Value = 47

$ws.Range("E34").Value = 16
$ws.Range("F34").Value = 8
$ws.Range("H34").Value = 8

$ws.Range("E36").Value = 92
$ws.Range("E38").Value = 64
$ws.Range("E41").Value = 33
$ws.Range("E48").Value = 27

$ws.Range("F66").Value = 18
$ws.Range("H66").Value = 18

$ws.Range("E70").Value = 38
$ws.Range("E71").Value = 29
$ws.Range("E77").Value = 50
$ws.Range("E81").Value = 14
$ws.Range("E84").Value = 4
